$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "RMSE"
$ws.Range("C1").Value = "NRMSE"
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = "RSE"
$ws.Range("F1").Value = "RRSE"
$ws.Range("G1").Value = "RAE"
$ws.Range("H1").Value = "R2"
$ws.Range("I1").Value = "Corr Coeff"

# Row 2 - random_forest
$ws.Range("A2").Value = "random_forest"
$ws.Range("B2").Value = 3.4294406323069455
$ws.Range("C2").Value = 0.24489724000178287
$ws.Range("D2").Value = 2.5701659599140245
$ws.Range("E2").Value = 0.26511098424749086
$ws.Range("F2").Value = 0.5148892931956256
$ws.Range("G2").Value = 0.47307042828657803
$ws.Range("H2").Value = 0.7348890157525092
$ws.Range("I2").Value = 0.85968621626487884

# Row 3 - lsboost
$ws.Range("A3").Value = "lsboost"
$ws.Range("B3").Value = 3.4441550996750858
$ws.Range("C3").Value = 0.24594800391138574
$ws.Range("D3").Value = 2.5948526922852282
$ws.Range("E3").Value = 0.26739085210487901
$ws.Range("F3").Value = 0.5170984936207792
$ws.Range("G3").Value = 0.47761432282023381
$ws.Range("H3").Value = 0.73260914789512099
$ws.Range("I3").Value = 0.85797076147681173

# Row 4 - neural_network
$ws.Range("A4").Value = "neural_network"
$ws.Range("B4").Value = 3.4831372157508604
$ws.Range("C4").Value = 0.24873172687379913
$ws.Range("D4").Value = 2.6095982674769935
$ws.Range("E4").Value = 0.2734779469461604
$ws.Range("F4").Value = 0.52295118983147981
$ws.Range("G4").Value = 0.48032842598714903
$ws.Range("H4").Value = 0.72652205305383966
$ws.Range("I4").Value = 0.85355535309139108
